$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several cards moved from the "Text Cards to implement:" column (C)
# over to the "Finished Cards:" column (D) now that they are working.

# "Council Room" is finished: clear it from C14, record it at D19
$ws.Range("C14").ClearContents()
$ws.Range("D19").Value = "Council Room"

# "Witch" is finished: clear it from C26, record it at D18
$ws.Range("C26").ClearContents()
$ws.Range("D18").Value = "Witch"

# "Moat" is finished: move it from C20 to D20
$ws.Range("C20").ClearContents()
$ws.Range("D20").Value = "Moat"

# Basic treasure/victory cards are finished: move them from column C to
# column D, each staying on the same row (28-33).
$ws.Range("C28").ClearContents()
$ws.Range("D28").Value = "Estate"

$ws.Range("C29").ClearContents()
$ws.Range("D29").Value = "Duchy"

$ws.Range("C30").ClearContents()
$ws.Range("D30").Value = "Province"

$ws.Range("C31").ClearContents()
$ws.Range("D31").Value = "Copper"

$ws.Range("C32").ClearContents()
$ws.Range("D32").Value = "Silver"

$ws.Range("C33").ClearContents()
$ws.Range("D33").Value = "Gold"

# Reflect the author's final selection/scroll position in the sheet.
$ws.Range("D20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
